# matchsync-ig asset refresh: rename "Matchsource"/"matchsource" -> "MatchSync"/"matchsync",
# bump the Date value, and record the (new) Experimental flag on the Metadata sheet.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Row 2 "URL": http://fhir.nmdp.org/ig/matchsource/... -> .../matchsync/...
$meta.Cells.Item(2, 2).Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/ms-rh-codes"

# Row 5 "Title": Matchsource Rh Value Sets -> MatchSync Rh Value Sets
$meta.Cells.Item(5, 2).Value = "MatchSync Rh Value Sets"

# Row 7 "Experimental": previously blank, now the literal text "true" (not boolean TRUE).
# A plain .Value assignment of "true" auto-coerces to a Boolean, so instead we write it
# as a quoted-string formula and immediately collapse the formula to its literal text
# result in place, which keeps the cell a text/shared-string cell.
$expCell = $meta.Cells.Item(7, 2)
$expCell.Formula = "=""true"""
$expCell.Copy()
$expCell.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Row 8 "Date": refreshed timestamp
$meta.Cells.Item(8, 2).Value = "2024-02-19T18:37:26-06:00"

# Row 11 "Description": MatchSource codes ... -> MatchSync codes ...
$meta.Cells.Item(11, 2).Value = "MatchSync codes for blood Rh. Combines NMDP and LOINC code"

# "Include ValueSets" sheet, row 2: matchsource -> matchsync in the NMDP status-codes URL
$incl = $wb.Worksheets.Item("Include ValueSets")
$incl.Cells.Item(2, 1).Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-rh-status-codes"
